# src/yarm/tests_data/test_validate_complete_config_valid/SOURCE_A.xlsx
#
# 1. Shared string "TOUCHED" -> "TOUCHED UP" (cell D1 on sheet "A.DATETIME").
# 2. Sheet "A.DATETIME" view: active cell / selection moves from E1 to D1,
#    and the view scrolls back so the top-left visible cell is A1 (was C1).
# 3. Default column width tweaks on both sheets (11.625 -> 11.66796875) —
#    applied best-effort via StandardWidth for completeness.

$wb = $excel.ActiveWorkbook

$wsDateTime = $wb.Worksheets.Item("A.DATETIME")
$wsA1       = $wb.Worksheets.Item("A.1")

# --- 1. Update the cell text ---------------------------------------------
$wsDateTime.Range("D1").Value = "TOUCHED UP"

# --- 3. Default column width (best effort) --------------------------------
$wsDateTime.StandardWidth = 11.66796875
$wsA1.StandardWidth = 11.66796875

# --- 2. Fix up the view / selection on the A.DATETIME sheet ---------------
$wsDateTime.Activate()

$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1

[void]$wsDateTime.Range("D1").Select()
